# Weekly fruit/vegetable price update: insert a new weekly record row
# before the existing row 51, shifting the old rows 51-52 down to 52-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51 (pushes current 51 -> 52, 52 -> 53)
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the latest weekly record
$ws.Range("A51").Value = 8
$ws.Range("B51").Value = "Terminal La Palmera de La Serena"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 45166
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 400
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 11000
$ws.Range("M51").Value = 10500
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Provincia del Elquí"
$ws.Range("P51").Value = 420
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
